# Hands-on tips: add a "設定加入好友訊息" (friend-add greeting) keyword row
# right after the header row, pushing the existing keyword rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (shifts rows 2-7 down to 3-8)
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new keyword + replies
$ws.Range("A2").Value = "設定加入好友訊息"
$ws.Range("B2").Value = "Hi"
$ws.Range("C2").Value = "HiHi "
$ws.Range("D2").Value = "HiHiHi"
$ws.Range("E2").Value = "HiHiHiHi"

# Match the final selection shown in the saved workbook
$ws.Range("D20").Select()
